# category_import_template.xlsx:
#  - rename first sheet "Dữ liệu mẫu" -> "Dữ liệu danh mục"
#  - make the second sheet ("Hướng dẫn") the active/selected tab
#    (was the first sheet before)

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "Dữ liệu danh mục"

$guideSheet = $wb.Worksheets.Item(2)
$guideSheet.Select()
$guideSheet.Activate()

$wb.Save()
